$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from H1 so I1/J1 match existing header formatting
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for columns I (I0) and J (IF)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 6
